# Applies the cell-value updates for the crypto price/volume refresh described
# in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
$updates = @(
    @('D2', '27.811.02'),
    @('E2', '  -0.67%  '),
    @('D3', '1.906.35'),
    @('E3', '  -0.03%  '),
    @('D4', '1.002'),
    @('E4', '  -0.25%  '),
    @('D5', '312.57'),
    @('E5', '  -1.61%  '),
    @('D6', '1.000'),
    @('E6', '  -0.32%  '),
    @('E7', '  +3.60%  '),
    @('D8', '0.3795'),
    @('E8', '  -0.10%  '),
    @('D9', '0.07269'),
    @('E9', '  -1.34%  '),
    @('B10', 'Polygon'),
    @('C10', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @('D10', '0.9074'),
    @('E10', '  -2.67%  '),
    @('B11', 'Solana'),
    @('C11', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'),
    @('D11', '21.19'),
    @('E11', '  +1.95%  '),
    @('B12', 'TRON'),
    @('C12', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'),
    @('D12', '0.07639'),
    @('E12', '  -1.40%  '),
    @('B13', 'WrappedEther'),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('D13', '1.900.87'),
    @('E13', '  -0.42%  '),
    @('D14', '5.462'),
    @('E14', '  -0.36%  '),
    @('D15', '92.22'),
    @('E15', '  +0.52%  '),
    @('D16', '1.003'),
    @('E16', '  -0.15%  '),
    @('D17', '0.000008718'),
    @('E17', '  -1.91%  '),
    @('D18', '0.9991'),
    @('D19', '27.853.19'),
    @('E19', '  -0.61%  '),
    @('D20', '14.61'),
    @('E20', '  -0.48%  '),
    @('D21', '5.169'),
    @('E21', '  +0.58%  '),
    @('D22', '2.161.53'),
    @('E22', '  +0.60%  '),
    @('D23', '10.85'),
    @('E23', '  -0.49%  '),
    @('D24', '6.592'),
    @('E24', '  -0.73%  '),
    @('D25', '152.71'),
    @('E25', '  -2.01%  '),
    @('D26', '1.842'),
    @('E26', '  -3.66%  '),
    @('D27', '2.216'),
    @('E27', '  +4.57%  '),
    @('D28', '18.36'),
    @('E28', '  -0.72%  '),
    @('D29', '114.92'),
    @('E29', '  -2.08%  '),
    @('D30', '4.882'),
    @('E30', '  -2.01%  '),
    @('D31', '0.08957'),
    @('E31', '  +0.13%  '),
    @('D32', '3.192'),
    @('E32', '  -1.75%  '),
    @('B33', 'Filecoin'),
    @('C33', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D33', '4.808'),
    @('E33', '  +3.02%  '),
    @('B34', 'ImmutableX'),
    @('C34', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D34', '0.7857'),
    @('E34', '  +1.94%  '),
    @('B35', 'ARBITRUM'),
    @('C35', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'),
    @('D35', '1.232'),
    @('E35', '  -1.58%  '),
    @('D36', '2.640'),
    @('E36', '  +1.73%  '),
    @('D37', '0.02078'),
    @('E37', '  +1.24%  '),
    @('D38', '3.059'),
    @('E38', '  +1.99%  '),
    @('D39', '1.092'),
    @('E39', '  -1.18%  '),
    @('D40', '0.5520'),
    @('E40', '  +0.11%  '),
    @('D41', '0.05276'),
    @('E41', '  -0.13%  '),
    @('D42', '6.774'),
    @('E42', '  -2.85%  '),
    @('D43', '113.38'),
    @('E43', '  +2.80%  '),
    @('D44', '8.464'),
    @('E44', '  -0.31%  '),
    @('E45', '  -1.23%  '),
    @('D46', '10.57'),
    @('E46', '  -0.73%  '),
    @('D47', '0.4804'),
    @('E47', '  -0.50%  '),
    @('D48', '0.9996'),
    @('E48', '  -0.38%  '),
    @('D49', '1.634'),
    @('E49', '  -0.70%  '),
    @('D50', '67.17'),
    @('E50', '  -1.02%  '),
    @('D51', '0.06043'),
    @('E51', '  -0.53%  ')
)

foreach ($update in $updates) {
    $addr = $update[0]
    $val = $update[1]
    $range = $ws.Range($addr)

    # Some "prices" are plain decimal-looking text (e.g. "1.002") that Excel would
    # otherwise auto-convert to a number. Force text entry with a leading quote,
    # then strip the formatting mark that the quote-prefix leaves behind so the
    # cell keeps the workbook default (unstyled) look, matching the source data,
    # which stores every value in this table as plain text.
    $isNumberLike = $val -match '^[+-]?\d+(\.\d+)?$'
    if ($isNumberLike) {
        $range.Value = "'" + $val
        $range.ClearFormats()
    } else {
        $range.Value = $val
    }
}
